$wb = $excel.ActiveWorkbook

$wsAbout = $wb.Worksheets.Item("About")
$wsFPIEBP = $wb.Worksheets.Item("FPIEBP")

# About sheet: update the "last updated" date in C1 from 1/3/2024 to 3/28/2024
$wsAbout.Range("C1").Value = "3/28/2024"

# FPIEBP sheet: hard coal row priorities changed (production/imports/exports)
# production: 3 -> 1, imports: 2 -> 3, exports: 1 -> 2
$wsFPIEBP.Range("B3").Value = 1
$wsFPIEBP.Range("C3").Value = 3
$wsFPIEBP.Range("D3").Value = 2

# Restore the active sheet and move the selection to E3 on FPIEBP
$wsFPIEBP.Activate() | Out-Null
$wsFPIEBP.Range("E3").Select() | Out-Null
